$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: point status -> Closed (resolution text already correct, untouched) ---
$ws.Range("F7").Value = "Closed"

$ws.Rows.Item(7).RowHeight = 60

# --- Row 8: fix typo in review point, clear acceptance, append reviewer note ---
$ws.Range("D8").Value = "Cur_x and Cur_y are not defined with specific values, developer shall know the exact values for them"

$ws.Range("E8").ClearContents()

$g8 = $ws.Range("G8")
$g8Existing = $g8.Characters().Text
$g8Text = $g8Existing + "`nMali 19/2/2020: I didn't mean that, I meant for example Req_PO1_DGC_SRS_009_V01 the requirement shall mention the value of x and y"
$g8.Value = $g8Text
foreach ($req in @("Req_PO1_DGC_SRS_014_V01", "Req_PO1_DGC_SRS_015_V01", "Req_PO1_DGC_SRS_016_V01", "Req_PO1_DGC_SRS_017_V01")) {
    $idx = $g8Text.IndexOf($req)
    $g8.Characters($idx + 1, $req.Length).Font.Bold = $true
}

$ws.Rows.Item(8).RowHeight = 165

# --- Restore the selection to reflect where the reviewer ended up ---
$ws.Range("G9").Select()
